$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 889129.6
$ws.Range("J19").Value = 257.22223
$ws.Range("L19").Value = 257.22223
$ws.Range("N19").Value = -607.2222300000001
# Row 28
$ws.Range("H28").Value = 729.46155
$ws.Range("I28").Value = 887.1
$ws.Range("K28").Value = 887.1
$ws.Range("M28").Value = -402.1
# Row 55
$ws.Range("H55").Value = 193.55556
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 220.28572
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 220.28572
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -648.28572
# Row 92
$ws.Range("H92").Value = 938
$ws.Range("I92").Value = 467.75
$ws.Range("J92").Value = 4700
$ws.Range("K92").Value = 467.75
$ws.Range("L92").Value = 4700
$ws.Range("M92").Value = 780.25
$ws.Range("N92").Value = -7196
# Row 112
$ws.Range("H112").Value = 2390.6
$ws.Range("J112").Value = 2589.1538
$ws.Range("L112").Value = 7767.4614
$ws.Range("N112").Value = -9983.4614
# Row 116
$ws.Range("H116").Value = 8660.5
$ws.Range("J116").Value = 11286.571
$ws.Range("L116").Value = 11286.571
$ws.Range("N116").Value = -18170.571
# Row 129
$ws.Range("H129").Value = 977.42426
$ws.Range("J129").Value = 992.375
$ws.Range("L129").Value = 2977.125
$ws.Range("N129").Value = -12977.125
# Row 132
$ws.Range("H132").Value = 45640216
$ws.Range("I132").Value = 52844830
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 158534490
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -158531960
$ws.Range("N132").Value = -38060
# Row 138
$ws.Range("H138").Value = 2901.68
$ws.Range("I138").Value = 1538.4546
$ws.Range("J138").Value = 3286.1794
$ws.Range("K138").Value = 4615.3638
$ws.Range("L138").Value = 9858.538199999999
$ws.Range("M138").Value = 524.6361999999999
$ws.Range("N138").Value = -20138.5382

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3581
$ws.Range("I2").Value = 902
$ws.Range("J2").Value = 7599.5
$ws.Range("K2").Value = 902
$ws.Range("L2").Value = 7599.5
$ws.Range("M2").Value = -789
$ws.Range("N2").Value = -7825.5
# Row 63
$ws.Range("H63").Value = 9897287
$ws.Range("I63").Value = 17316064
$ws.Range("J63").Value = 5583.3335
$ws.Range("K63").Value = 17316064
$ws.Range("L63").Value = 5583.3335
$ws.Range("M63").Value = -17315378
$ws.Range("N63").Value = -6955.3335
# Row 66
$ws.Range("H66").Value = 9897287
$ws.Range("I66").Value = 17316064
$ws.Range("J66").Value = 5583.3335
$ws.Range("K66").Value = 86580320
$ws.Range("L66").Value = 27916.6675
$ws.Range("M66").Value = -86576888
$ws.Range("N66").Value = -34780.6675
# Row 88
$ws.Range("H88").Value = 13335013
$ws.Range("I88").Value = 16668141
$ws.Range("K88").Value = 16668141
$ws.Range("M88").Value = -16667735
# Row 91
$ws.Range("H91").Value = 13335013
$ws.Range("I91").Value = 16668141
$ws.Range("K91").Value = 16668141
$ws.Range("M91").Value = -16666737
# Row 116
$ws.Range("H116").Value = 3581
$ws.Range("I116").Value = 902
$ws.Range("J116").Value = 7599.5
$ws.Range("K116").Value = 902
$ws.Range("L116").Value = 7599.5
$ws.Range("M116").Value = 1392
$ws.Range("N116").Value = -12187.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3581
$ws.Range("I3").Value = 902
$ws.Range("J3").Value = 7599.5
$ws.Range("K3").Value = 902
$ws.Range("L3").Value = 7599.5
$ws.Range("M3").Value = -788
$ws.Range("N3").Value = -7827.5
# Row 86
$ws.Range("H86").Value = 2657.1428
$ws.Range("I86").Value = 2025
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 2025
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -902
$ws.Range("N86").Value = -5746
# Row 89
$ws.Range("H89").Value = 2657.1428
$ws.Range("I89").Value = 2025
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 10125
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -4509
$ws.Range("N89").Value = -28732
# Row 105
$ws.Range("H105").Value = 2211.5264
$ws.Range("I105").Value = 2116.923
$ws.Range("J105").Value = 2416.5
$ws.Range("K105").Value = 2116.923
$ws.Range("L105").Value = 2416.5
$ws.Range("M105").Value = -369.9229999999998
$ws.Range("N105").Value = -5910.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3423.8044
$ws.Range("I31").Value = 1397.7931
$ws.Range("J31").Value = 6879.9414
$ws.Range("K31").Value = 1397.7931
$ws.Range("L31").Value = 6879.9414
$ws.Range("M31").Value = -1102.7931
$ws.Range("N31").Value = -7469.9414
# Row 32
$ws.Range("H32").Value = 11500
$ws.Range("I32").Value = 11500
$ws.Range("K32").Value = 11500
$ws.Range("M32").Value = -11184
# Row 34
$ws.Range("H34").Value = 3423.8044
$ws.Range("I34").Value = 1397.7931
$ws.Range("J34").Value = 6879.9414
$ws.Range("K34").Value = 1397.7931
$ws.Range("L34").Value = 6879.9414
$ws.Range("M34").Value = -1195.7931
$ws.Range("N34").Value = -7283.9414
# Row 38
$ws.Range("H38").Value = 49999
$ws.Range("J38").Value = 49999
$ws.Range("L38").Value = 49999
$ws.Range("N38").Value = -50753
# Row 39
$ws.Range("H39").Value = 15127.786
$ws.Range("I39").Value = 2689.8572
$ws.Range("J39").Value = 27565.715
$ws.Range("K39").Value = 2689.8572
$ws.Range("L39").Value = 27565.715
$ws.Range("M39").Value = -2298.8572
$ws.Range("N39").Value = -28347.715
# Row 46
$ws.Range("H46").Value = 49999
$ws.Range("J46").Value = 49999
$ws.Range("L46").Value = 49999
$ws.Range("N46").Value = -50421
# Row 49
$ws.Range("H49").Value = 15127.786
$ws.Range("I49").Value = 2689.8572
$ws.Range("J49").Value = 27565.715
$ws.Range("K49").Value = 2689.8572
$ws.Range("L49").Value = 27565.715
$ws.Range("M49").Value = -2507.8572
$ws.Range("N49").Value = -27929.715
# Row 62
$ws.Range("H62").Value = 33338730
$ws.Range("I62").Value = 71433520
$ws.Range("J62").Value = 5791.75
$ws.Range("K62").Value = 71433520
$ws.Range("L62").Value = 5791.75
$ws.Range("M62").Value = -71432896
$ws.Range("N62").Value = -7039.75
# Row 65
$ws.Range("H65").Value = 33338730
$ws.Range("I65").Value = 71433520
$ws.Range("J65").Value = 5791.75
$ws.Range("K65").Value = 357167600
$ws.Range("L65").Value = 28958.75
$ws.Range("M65").Value = -357164480
$ws.Range("N65").Value = -35198.75
# Row 99
$ws.Range("H99").Value = 3770.611
$ws.Range("I99").Value = 2270.3635
$ws.Range("K99").Value = 2270.3635
$ws.Range("M99").Value = -772.3634999999999
# Row 126
$ws.Range("H126").Value = 3770.611
$ws.Range("I126").Value = 2270.3635
$ws.Range("K126").Value = 6811.0905
$ws.Range("M126").Value = -4341.0905
# Row 141
$ws.Range("H141").Value = 19154.055
$ws.Range("J141").Value = 19154.055
$ws.Range("L141").Value = 19154.055
$ws.Range("N141").Value = -29514.055

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 661.17645
$ws.Range("J113").Value = 618.625
$ws.Range("L113").Value = 1855.875
$ws.Range("N113").Value = -6195.875
# Row 121
$ws.Range("H121").Value = 1223.6451
$ws.Range("I121").Value = 406.66666
$ws.Range("J121").Value = 1265.1864
$ws.Range("K121").Value = 1219.99998
$ws.Range("L121").Value = 3795.5592
$ws.Range("M121").Value = 90.00001999999995
$ws.Range("N121").Value = -6415.5592
# Row 129
$ws.Range("H129").Value = 2562.6
$ws.Range("I129").Value = 2522
$ws.Range("J129").Value = 2593.05
$ws.Range("K129").Value = 7566
$ws.Range("L129").Value = 7779.150000000001
$ws.Range("M129").Value = -2566
$ws.Range("N129").Value = -17779.15
# Row 137
$ws.Range("H137").Value = 2839.3333
$ws.Range("I137").Value = 2829.2856
$ws.Range("J137").Value = 2980
$ws.Range("K137").Value = 8487.856800000001
$ws.Range("L137").Value = 8940
$ws.Range("M137").Value = -3387.856800000001
$ws.Range("N137").Value = -19140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5589.8887
$ws.Range("I70").Value = 5468.9443
$ws.Range("K70").Value = 5468.9443
$ws.Range("M70").Value = -5198.9443
# Row 73
$ws.Range("H73").Value = 5589.8887
$ws.Range("I73").Value = 5468.9443
$ws.Range("K73").Value = 5468.9443
$ws.Range("M73").Value = -4532.9443
# Row 80
$ws.Range("H80").Value = 25002560
$ws.Range("I80").Value = 50002120
$ws.Range("K80").Value = 50002120
$ws.Range("M80").Value = -50001122
# Row 83
$ws.Range("H83").Value = 25002560
$ws.Range("I83").Value = 50002120
$ws.Range("K83").Value = 250010600
$ws.Range("M83").Value = -250005608
# Row 141
$ws.Range("H141").Value = 63107.5
$ws.Range("J141").Value = 67476.664
$ws.Range("L141").Value = 67476.664
$ws.Range("N141").Value = -77836.664

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 2262.5
$ws.Range("J2").Value = 2262.5
$ws.Range("L2").Value = 2262.5
$ws.Range("N2").Value = -2486.5
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 45
$ws.Range("H45").Value = 29000
$ws.Range("I45").Value = 29000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 29000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -28593
$ws.Range("N45").ClearContents()
# Row 122
$ws.Range("H122").Value = 6255
$ws.Range("I122").Value = 3040
$ws.Range("K122").Value = 9120
$ws.Range("M122").Value = -6670

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1315.3334
$ws.Range("J107").Value = 1555.6
$ws.Range("L107").Value = 4666.799999999999
$ws.Range("N107").Value = -8506.799999999999
# Row 126
$ws.Range("H126").Value = 274709.72
$ws.Range("I126").Value = 1095.1923
$ws.Range("J126").Value = 821938.75
$ws.Range("K126").Value = 3285.5769
$ws.Range("L126").Value = 2465816.25
$ws.Range("M126").Value = -815.5769
$ws.Range("N126").Value = -2470756.25
# Row 135
$ws.Range("H135").Value = 79044.086
$ws.Range("J135").Value = 79044.086
$ws.Range("L135").Value = 79044.086
$ws.Range("N135").Value = -89184.086
# Row 140
$ws.Range("H140").Value = 51229.125
$ws.Range("J140").Value = 51229.125
$ws.Range("L140").Value = 51229.125
$ws.Range("N140").Value = -61589.125
# Row 141
$ws.Range("H141").Value = 44000
$ws.Range("J141").Value = 44000
$ws.Range("L141").Value = 44000
$ws.Range("N141").Value = -54360
